$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update Quantite column (C) values to reflect stock adjustments
$ws.Range("C2").Value = 29
$ws.Range("C3").Value = 29
$ws.Range("C5").Value = 1226
$ws.Range("C7").Value = 7
